$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value. Cells whose new text would otherwise be auto-coerced
# into a Number by Excel (losing a trailing zero, e.g. "11.40" -> 11.4)
# are written with an explicit text NumberFormat, then restored to the
# default "Normal" style so no stray formatting is left behind.
$updates = @(
    @{ Cell = 'D2'; Value = '34.064.88'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -0.18%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.789.04'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -0.08%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  +0.05%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '226.94'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +1.80%  '; ForceText = $false }
    @{ Cell = 'E6'; Value = '  -1.33%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  +0.04%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '32.34'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  -0.11%  '; ForceText = $false }
    @{ Cell = 'E9'; Value = '  +3.98%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.0685'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -4.38%  '; ForceText = $false }
    @{ Cell = 'E11'; Value = '  +1.00%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '2.046.61'; ForceText = $false }
    @{ Cell = 'E12'; Value = '  -0.02%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '11.40'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  +3.97%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '1.824.75'; ForceText = $false }
    @{ Cell = 'E14'; Value = '  +2.17%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '0.623'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  -0.51%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '34.053.26'; ForceText = $false }
    @{ Cell = 'E16'; Value = '  -0.08%  '; ForceText = $false }
    @{ Cell = 'E17'; Value = '  +0.47%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '67.94'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  -0.21%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '242.46'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  -0.91%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '0.0₃0775'; ForceText = $false }
    @{ Cell = 'E20'; Value = '  -1.21%  '; ForceText = $false }
    @{ Cell = 'E21'; Value = '  +0.02%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '10.74'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -0.06%  '; ForceText = $false }
    @{ Cell = 'E23'; Value = '  +0.03%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '162.09'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  +2.06%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '7.15'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  +1.22%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '16.22'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  -0.99%  '; ForceText = $false }
    @{ Cell = 'E28'; Value = '  +0.63%  '; ForceText = $false }
    @{ Cell = 'E29'; Value = '  +0.13%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '1.23'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  +2.30%  '; ForceText = $false }
    @{ Cell = 'E31'; Value = '  -0.76%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '3.65'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -0.79%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '3.62'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  +3.34%  '; ForceText = $false }
    @{ Cell = 'E34'; Value = '  +1.85%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '1.398.31'; ForceText = $false }
    @{ Cell = 'E35'; Value = '  +0.36%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '0.654'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  +0.67%  '; ForceText = $false }
    @{ Cell = 'B37'; Value = 'RenderToken'; ForceText = $false }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; ForceText = $false }
    @{ Cell = 'D37'; Value = '2.37'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  +9.07%  '; ForceText = $false }
    @{ Cell = 'B38'; Value = 'TrustWalletToken'; ForceText = $false }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; ForceText = $false }
    @{ Cell = 'D38'; Value = '1.04'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  -0.77%  '; ForceText = $false }
    @{ Cell = 'E39'; Value = '  +1.44%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '80.15'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  +0.37%  '; ForceText = $false }
    @{ Cell = 'E41'; Value = '  -0.02%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '0.920'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  +0.00%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '13.73'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  +14.29%  '; ForceText = $false }
    @{ Cell = 'E44'; Value = '  -0.41%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '6.11'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  +2.45%  '; ForceText = $false }
    @{ Cell = 'E46'; Value = '  +8.28%  '; ForceText = $false }
    @{ Cell = 'B47'; Value = 'Kaspa'; ForceText = $false }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; ForceText = $false }
    @{ Cell = 'D47'; Value = '0.0506'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  +1.85%  '; ForceText = $false }
    @{ Cell = 'B48'; Value = 'WEMIXToken'; ForceText = $false }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; ForceText = $false }
    @{ Cell = 'D48'; Value = '1.08'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  +2.49%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '107.62'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -0.03%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '1.947.53'; ForceText = $false }
    @{ Cell = 'E50'; Value = '  -0.12%  '; ForceText = $false }
    @{ Cell = 'E51'; Value = '  +0.05%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = '@'
        $rng.Value = $u.Value
        $rng.Style = 'Normal'
    } else {
        $rng.Value = $u.Value
    }
}
